$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Plan")

# --- Row 19: Art: Animations -> due 24-Nov-2021, Status: Not Implemented ---
$ws.Range("B19").Value = 44524
$ws.Range("B19").NumberFormat = "d-mmm"
$ws.Range("B19").HorizontalAlignment = -4108
$ws.Range("B19").VerticalAlignment = -4108
$ws.Range("B19").WrapText = $true
$ws.Range("D19").Value = "Not Implemented"

# --- Row 20: UI: On Screen Text -> "Ongoing as needed", Status: In Progress ---
$ws.Range("B20").Value = "Ongoing as needed"
$ws.Range("D20").Value = "In Progress"

# --- Row 21: Art: Tileset -> due 17-Nov-2021, Status: In Progress ---
$ws.Range("B21").Value = 44517
$ws.Range("B21").NumberFormat = "d-mmm"
$ws.Range("B21").HorizontalAlignment = -4108
$ws.Range("B21").VerticalAlignment = -4108
$ws.Range("B21").WrapText = $true
$ws.Range("D21").Value = "In Progress"

# --- Row 22: Art: Character -> due 17-Nov-2021, Status: In Progress ---
$ws.Range("B22").Value = 44517
$ws.Range("B22").NumberFormat = "d-mmm"
$ws.Range("B22").HorizontalAlignment = -4108
$ws.Range("B22").VerticalAlignment = -4108
$ws.Range("B22").WrapText = $true
$ws.Range("D22").Value = "In Progress"

# --- Row 23: Art: Enemies -> due 24-Nov-2021, Status: Not Implemented ---
$ws.Range("B23").Value = 44524
$ws.Range("B23").NumberFormat = "d-mmm"
$ws.Range("B23").HorizontalAlignment = -4108
$ws.Range("B23").VerticalAlignment = -4108
$ws.Range("B23").WrapText = $true
$ws.Range("D23").Value = "Not Implemented"

# --- View state: scroll window back to top-left, move selection to B23 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
